# Notification against transaction implement Queues
#
# Adds a "Current User" recipient to several notification rows and a new
# "Api" status/queue column (F) tracking API implementation progress on the
# "Notifications" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notifications")

# --- Row 2: On Submit Challenge -> add F2 "Api" ---
$ws.Range("F2").Value = "Api"
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4108

# --- Row 3: On First Vote -> add F3 "Api" ---
$ws.Range("F3").Value = "Api"
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F3").VerticalAlignment = -4108

# --- Row 4: On Second Vote -> unchanged ---

# --- Row 5: On Load Balance -> recipient gains Current User; Done/Api added ---
$ws.Range("B5").Value = "to Admin, Current User"
$ws.Range("E5").Value = "Done"
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("E5").VerticalAlignment = -4108
$ws.Range("F5").Value = "Api"
$ws.Range("F5").HorizontalAlignment = -4108
$ws.Range("F5").VerticalAlignment = -4108

# --- Row 6: On Miscellaneous -> recipient gains Current User; Done/Api added ---
$ws.Range("B6").Value = "to Admin, Current User"
$ws.Range("E6").Value = "Done"
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E6").VerticalAlignment = -4108
$ws.Range("F6").Value = "Api"
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("F6").VerticalAlignment = -4108

# --- Row 7: On Withdraw -> recipient gains Current User; Done/Api added ---
$ws.Range("B7").Value = "to Admin, Current User"
$ws.Range("E7").Value = "Done"
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").VerticalAlignment = -4108
$ws.Range("F7").Value = "Api"
$ws.Range("F7").HorizontalAlignment = -4108
$ws.Range("F7").VerticalAlignment = -4108

# --- Row 8: On Donate -> recipient gains Current User; status ++ -> Done; Api added ---
$ws.Range("B8").Value = "to Admin, Creater, Current User"
$ws.Range("E8").Value = "Done"
$ws.Range("F8").Value = "Api"
$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("F8").VerticalAlignment = -4108

# --- Row 9: On Create Challenge -> recipient gains Current User ---
$ws.Range("B9").Value = "to Admin, Current User"

# --- Row 10: On Win -> unchanged ---

# Update the active selection to reflect the new working cell (G8).
$ws.Range("G8").Select()
